$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("R/RStudio", $true, $false, $false, $false, $false, $true, 1, $false, "R, RStudio", 2)
